# Swap the data values between row 2 and row 3 for the columns that
# actually differ between the two records: A, B, D, E, F, G, H, Q, R.
# (Other columns happen to hold identical values in both rows, so a
# full-row swap and this targeted swap are equivalent for this sheet.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$columns = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

foreach ($col in $columns) {
    $cellRow2 = $ws.Range("$col" + "2")
    $cellRow3 = $ws.Range("$col" + "3")

    $value2 = $cellRow2.Value()
    $value3 = $cellRow3.Value()

    $cellRow2.Value = $value3
    $cellRow3.Value = $value2
}
